# Add a new "Questions?" slide to the end of the presentation.
#
# The original commit ("Add questions slides to end of ppts") appends a
# single new slide after the last existing slide ("Helpful Tips"). The new
# slide uses the "Animated Closing Slide" custom layout (the layout that
# exposes a single centered `body`/idx=10 placeholder, normally prompting
# "<Call to action>"), and simply fills that placeholder with "Questions?".

$p = $ppt.ActivePresentation

# The "Animated Closing Slide" layout is the 32nd custom layout on the
# slide master; it contributes exactly one placeholder shape to new
# slides ("Text Placeholder 1", type="body" sz="quarter" idx="10").
$layout = $p.SlideMaster.CustomLayouts.Item(32)

# Append as the new last slide.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.AddSlide($newIndex, $layout)

# Fill in the placeholder text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"

# Match the fade transition used throughout the rest of the deck.
$s.SlideShowTransition.EntryEffect = 3849
